$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 363.42856
$ws.Range("I33").Value = 350.77777
$ws.Range("K33").Value = 350.77777
$ws.Range("M33").Value = -121.77777
$ws.Range("H74").Value = 6479
$ws.Range("I74").Value = 6727.143
$ws.Range("J74").Value = 5900
$ws.Range("K74").Value = 6727.143
$ws.Range("L74").Value = 5900
$ws.Range("M74").Value = -5791.143
$ws.Range("N74").Value = -7772
$ws.Range("H77").Value = 6479
$ws.Range("I77").Value = 6727.143
$ws.Range("J77").Value = 5900
$ws.Range("K77").Value = 33635.715
$ws.Range("L77").Value = 29500
$ws.Range("M77").Value = -28955.715
$ws.Range("N77").Value = -38860
$ws.Range("H92").Value = 699.8570999999999
$ws.Range("I92").Value = 757.55554
$ws.Range("J92").Value = 596
$ws.Range("K92").Value = 757.55554
$ws.Range("L92").Value = 596
$ws.Range("M92").Value = 490.44446
$ws.Range("N92").Value = -3092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1432.75
$ws.Range("I2").Value = 910.5
$ws.Range("K2").Value = 910.5
$ws.Range("M2").Value = -797.5
$ws.Range("H32").Value = 3239.88
$ws.Range("I32").Value = 2541.5417
$ws.Range("K32").Value = 2541.5417
$ws.Range("M32").Value = -2254.5417
$ws.Range("H61").Value = 1829.6
$ws.Range("I61").Value = 1712
$ws.Range("K61").Value = 1712
$ws.Range("M61").Value = -1500
$ws.Range("H88").Value = 1407.1666
$ws.Range("I88").Value = 1046.8
$ws.Range("K88").Value = 1046.8
$ws.Range("M88").Value = -640.8
$ws.Range("H91").Value = 1407.1666
$ws.Range("I91").Value = 1046.8
$ws.Range("K91").Value = 1046.8
$ws.Range("M91").Value = 357.2
$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344
$ws.Range("H116").Value = 1432.75
$ws.Range("I116").Value = 910.5
$ws.Range("K116").Value = 910.5
$ws.Range("M116").Value = 1383.5
$ws.Range("H122").Value = 5268
$ws.Range("I122").Value = 5268
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15804
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -13354
$ws.Range("H131").Value = 59626.668
$ws.Range("J131").Value = 59626.668
$ws.Range("L131").Value = 59626.668
$ws.Range("N131").Value = -69706.66800000001
$ws.Range("H132").Value = 5232.3076
$ws.Range("I132").Value = 3791.6316
$ws.Range("J132").Value = 9142.714
$ws.Range("K132").Value = 11374.8948
$ws.Range("L132").Value = 27428.142
$ws.Range("M132").Value = -8844.8948
$ws.Range("N132").Value = -32488.142
$ws.Range("H136").Value = 1829.6
$ws.Range("I136").Value = 1712
$ws.Range("K136").Value = 5136
$ws.Range("M136").Value = -2586

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1432.75
$ws.Range("I3").Value = 910.5
$ws.Range("K3").Value = 910.5
$ws.Range("M3").Value = -796.5
$ws.Range("H20").Value = 1983.0834
$ws.Range("I20").Value = 2199.7
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 2199.7
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = -1952.7
$ws.Range("N20").Value = -1394
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null
$ws.Range("H86").Value = 9806.706
$ws.Range("J86").Value = 3896.7778
$ws.Range("L86").Value = 3896.7778
$ws.Range("N86").Value = -6142.7778
$ws.Range("H89").Value = 9806.706
$ws.Range("J89").Value = 3896.7778
$ws.Range("L89").Value = 19483.889
$ws.Range("N89").Value = -30715.889
$ws.Range("H134").Value = 1999
$ws.Range("I134").Value = 1999
$ws.Range("K134").Value = 5997
$ws.Range("M134").Value = -3462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4255.5
$ws.Range("I31").Value = 5948
$ws.Range("K31").Value = 5948
$ws.Range("M31").Value = -5653
$ws.Range("H34").Value = 4255.5
$ws.Range("I34").Value = 5948
$ws.Range("K34").Value = 5948
$ws.Range("M34").Value = -5746
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256
$ws.Range("H99").Value = 9128
$ws.Range("I99").Value = 9128
$ws.Range("K99").Value = 9128
$ws.Range("M99").Value = -7630
$ws.Range("H126").Value = 9128
$ws.Range("I126").Value = 9128
$ws.Range("K126").Value = 27384
$ws.Range("M126").Value = -24914
$ws.Range("H132").Value = 4348.077
$ws.Range("J132").Value = 4721.3335
$ws.Range("L132").Value = 14164.0005
$ws.Range("N132").Value = -19224.0005
$ws.Range("H134").Value = 1601.25
$ws.Range("I134").Value = 1410.5454
$ws.Range("K134").Value = 4231.6362
$ws.Range("M134").Value = -1696.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6409.25
$ws.Range("J34").Value = 6409.25
$ws.Range("L34").Value = 19227.75
$ws.Range("N34").Value = -19395.75
$ws.Range("H121").Value = 5333.1665
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 5999.8
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 17999.4
$ws.Range("M121").Value = -4690
$ws.Range("N121").Value = -20619.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9937.700000000001
$ws.Range("I70").Value = 9497.444
$ws.Range("K70").Value = 9497.444
$ws.Range("M70").Value = -9227.444
$ws.Range("H73").Value = 9937.700000000001
$ws.Range("I73").Value = 9497.444
$ws.Range("K73").Value = 9497.444
$ws.Range("M73").Value = -8561.444
$ws.Range("H102").Value = 1052
$ws.Range("I102").Value = 1052
$ws.Range("K102").Value = 1052
$ws.Range("M102").Value = 570
$ws.Range("H122").Value = 2176.348
$ws.Range("I122").Value = 2176.348
$ws.Range("K122").Value = 6529.044
$ws.Range("M122").Value = -4079.044
$ws.Range("H126").Value = 3249.25
$ws.Range("I126").Value = 2498.5
$ws.Range("K126").Value = 7495.5
$ws.Range("M126").Value = -5025.5
$ws.Range("H128").Value = 94800
$ws.Range("I128").Value = 94800
$ws.Range("K128").Value = 94800
$ws.Range("M128").Value = -89820
$ws.Range("H132").Value = 2968.3447
$ws.Range("I132").Value = 2927
$ws.Range("J132").Value = 3226.75
$ws.Range("K132").Value = 8781
$ws.Range("L132").Value = 9680.25
$ws.Range("M132").Value = -6251
$ws.Range("N132").Value = -14740.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 7800.6
$ws.Range("I21").Value = 6999
$ws.Range("J21").Value = 11007
$ws.Range("K21").Value = 6999
$ws.Range("L21").Value = 11007
$ws.Range("M21").Value = -6825
$ws.Range("N21").Value = -11355
$ws.Range("H22").Value = 2719.6
$ws.Range("I22").Value = 849.5
$ws.Range("J22").Value = 3966.3333
$ws.Range("K22").Value = 849.5
$ws.Range("L22").Value = 3966.3333
$ws.Range("M22").Value = -554.5
$ws.Range("N22").Value = -4556.3333
$ws.Range("H27").Value = 2719.6
$ws.Range("I27").Value = 849.5
$ws.Range("J27").Value = 3966.3333
$ws.Range("K27").Value = 849.5
$ws.Range("L27").Value = 3966.3333
$ws.Range("M27").Value = -742.5
$ws.Range("N27").Value = -4180.3333
$ws.Range("H40").Value = 1726.0714
$ws.Range("I40").Value = 1388.75
$ws.Range("K40").Value = 1388.75
$ws.Range("M40").Value = -1252.75
$ws.Range("H122").Value = 2593.7856
$ws.Range("I122").Value = 2442.8333
$ws.Range("K122").Value = 7328.499899999999
$ws.Range("M122").Value = -4878.499899999999
$ws.Range("H128").Value = 55996
$ws.Range("J128").Value = 55996
$ws.Range("L128").Value = 55996
$ws.Range("N128").Value = -65956
$ws.Range("H130").Value = 62997
$ws.Range("J130").Value = 62997
$ws.Range("L130").Value = 62997
$ws.Range("N130").Value = -73037

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2355.1428
$ws.Range("I81").Value = 2205.8333
$ws.Range("K81").Value = 4411.6666
$ws.Range("M81").Value = -3350.6666
$ws.Range("H82").Value = 30301
$ws.Range("J82").Value = 30301
$ws.Range("L82").Value = 30301
$ws.Range("N82").Value = -31067
$ws.Range("H84").Value = 2355.1428
$ws.Range("I84").Value = 2205.8333
$ws.Range("K84").Value = 22058.333
$ws.Range("M84").Value = -16754.333
$ws.Range("H85").Value = 30301
$ws.Range("J85").Value = 30301
$ws.Range("L85").Value = 30301
$ws.Range("N85").Value = -32953
$ws.Range("H126").Value = 9800.6
$ws.Range("I126").Value = 9800.6
$ws.Range("K126").Value = 29401.8
$ws.Range("M126").Value = -26931.8
$ws.Range("H132").Value = 2698.375
$ws.Range("I132").Value = 2698.375
$ws.Range("K132").Value = 8095.125
$ws.Range("M132").Value = -5565.125
